$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("the red pill")

# Delete the first row (stray duplicate header), shifting all data up by one.
# This turns A1:H27 (header + 26 data rows) into A1:H26 (26 data rows).
$ws.Rows.Item(1).Delete()

# The named range that pointed at the data block needs to track the shift:
# it used to start at row 3 / end at row 28, now starts at row 2 / ends at row 27.
$wb.Names.Item("reddit_theredpill").RefersTo = '=''the red pill''!$A$2:$H$27'

# Refresh the sort state saved on the sheet so it reflects the new data extent.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A1"))
$ws.Sort.SetRange($ws.Range("A1:H27"))
$ws.Sort.Header = 0
$ws.Sort.Apply()
